$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the MDR introduction start-time parameter to the new unified name
$ws.Range("A5").Value = "mdr_introduce_time"

# Remove the now-obsolete "end_mdr_introduce_time" row entirely (old row 6),
# shifting the rows below it up by one
$ws.Rows.Item(6).Delete()

# Update the active cell selection to match the edited sheet
$ws.Range("B6").Select()
